$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-18 14:18:41"
$wsZh.Range("E4").Value = "2016-03-18 14:18:41"
$wsZh.Range("H3").Value = "2016-03-18 14:19:00"
$wsZh.Range("H4").Value = "2016-03-18 14:19:00"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-18 14:18:44"
$wsDe.Range("E4").Value = "2016-03-18 14:18:44"
$wsDe.Range("H3").Value = "2016-03-18 14:19:06"
$wsDe.Range("H4").Value = "2016-03-18 14:19:06"
